$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 93 data appended after the existing 02/24/2026 row.
# Column A holds a date-formatted text label (matches the existing rows,
# which store the date as literal text, not an Excel date serial), so we
# temporarily force a text number format, assign the value, then clear the
# formatting back to the sheet's default (General / no explicit style) -
# this keeps the stored cell type as Text without leaving a stray style
# behind on the cell.
$ws.Range("A93").NumberFormat = "@"
$ws.Range("A93").Value = "02/25/2026"
$ws.Range("A93").ClearFormats()

$ws.Range("B93").Value = 9167.34
$ws.Range("C93").Value = 0.2470002370263118
$ws.Range("D93").Value = 0.7529997629736882
$ws.Range("E93").Value = -341.23
$ws.Range("F93").Value = -36.5
$ws.Range("G93").Value = -24088.75
$ws.Range("H93").Value = -77.73
$ws.Range("I93").Value = -1188.52
$ws.Range("J93").Value = -34.42
$ws.Range("K93").Value = -25277.27
$ws.Range("L93").Value = -73.39
